# Update cryptos list with latest scraped values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "37.276.76"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.007.04"
$ws.Range("E3").Value = "  -1.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'259.44"
$ws.Range("E5").Value = "  +4.41%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.618"
$ws.Range("E6").Value = "  -1.64%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - Solana
$ws.Range("D8").Value = "'56.58"
$ws.Range("E8").Value = "  -6.78%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -3.39%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.0772"
$ws.Range("E10").Value = "  -4.82%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -3.17%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.303.32"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "'14.25"
$ws.Range("E13").Value = "  -6.93%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "'21.70"
$ws.Range("E14").Value = "  -2.94%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -7.78%  "

# Row 16 - Polkadot
$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  -5.98%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.029.35"
$ws.Range("E17").Value = "  -0.37%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "37.238.70"

# Row 19 - Litecoin
$ws.Range("D19").Value = "'70.18"
$ws.Range("E19").Value = "  -0.81%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -3.83%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'232.69"
$ws.Range("E21").Value = "  +0.65%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.07%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.08%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -1.28%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.71%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'164.73"
$ws.Range("E26").Value = "  +0.68%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'8.97"
$ws.Range("E27").Value = "  -5.43%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'19.57"
$ws.Range("E28").Value = "  -1.33%  "

# Row 29 - Kaspa
$ws.Range("E29").Value = "  -6.27%  "

# Row 30 - ImmutableX
$ws.Range("E30").Value = "  -4.41%  "

# Row 31 - Stellar
$ws.Range("D31").Value = "'0.119"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'4.59"
$ws.Range("E32").Value = "  -5.05%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -5.58%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").Value = "'4.45"
$ws.Range("E34").Value = "  -1.38%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  -5.97%  "

# Row 36 - RenderToken (was WEMIXToken)
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'3.40"
$ws.Range("E36").Value = "  -6.44%  "

# Row 37 - WEMIXToken (was BinanceUSD)
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +0.23%  "

# Row 38 - BinanceUSD (was RenderToken)
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  -0.15%  "

# Row 39 - THORChain
$ws.Range("D39").Value = "'5.45"
$ws.Range("E39").Value = "  -0.51%  "

# Row 40 - HuobiToken
$ws.Range("E40").Value = "  +1.79%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  -1.15%  "

# Row 42 - VeChain
$ws.Range("E42").Value = "  -1.79%  "

# Row 43 - Cronos
$ws.Range("E43").Value = "  -5.84%  "

# Row 44 - Maker
$ws.Range("D44").Value = "1.429.68"
$ws.Range("E44").Value = "  +3.12%  "

# Row 45 - Aave
$ws.Range("D45").Value = "'89.31"
$ws.Range("E45").Value = "  -3.69%  "

# Row 46 - InjectiveProtocol
$ws.Range("E46").Value = "  -9.44%  "

# Row 47 - ARBITRUM
$ws.Range("E47").Value = "  -3.87%  "

# Row 48 - MXToken
$ws.Range("E48").Value = "  +2.13%  "

# Row 49 - FraxShare
$ws.Range("D49").Value = "'6.97"
$ws.Range("E49").Value = "  -6.92%  "

# Row 50 - RocketPoolETH
$ws.Range("D50").Value = "2.195.40"
$ws.Range("E50").Value = "  -1.25%  "

# Row 51 - NEARProtocol
$ws.Range("D51").Value = "'1.95"
$ws.Range("E51").Value = "  -10.26%  "

